$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace the old sample email with a fresh one ---------------
# Clearing B2 first purges the now-unreferenced "sebastien.debeauffort@..."
# shared string so the remaining strings ("email"/"names") compact down,
# matching how the real edit history re-keyed the shared-strings table.
$ws.Range("B2").ClearContents()
$ws.Range("A2").Value = "Jean exemple"
$ws.Range("B2").Value = "jeanexemple@outlook.com"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:jeanexemple@outlook.com")

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Value = "Jean exemple2"
$ws.Range("B3").Value = "jeanexemple2@outlook.com"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:jeanexemple2@outlook.com")

# --- Row 4 ---------------------------------------------------------------
$ws.Range("A4").Value = "Jean exemple3"
$ws.Range("B4").Value = "jeanéxemple3@outlook.com"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:jeanéxemple3@outlook.com")

# --- Row 5 (hyperlink target / display text drift from the visible text) -
$ws.Range("A5").Value = "Jean exemple4"
$ws.Range("B5").Value = "jeanexemple2@outlook.com"
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:jeanexemple2@outlook.com", "", "", "jeanexemple2@outlook.com")
$ws.Range("B5").Value = "jeanexemple 4@outlook.com"

# --- Row 6 (hyperlink target / display text drift from the visible text) -
$ws.Range("A6").Value = "Jeanexemple5"
$ws.Range("B6").Value = "jeanexemple5@outlook.com"
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:jeanexemple5@outlook.com", "", "", "jeanexemple5@outlook.com")
$ws.Range("B6").Value = "jeanexemple5outlook.com"

# --- Row 7 ---------------------------------------------------------------
$ws.Range("A7").Value = "Jeanexemple6"
$ws.Range("B7").Value = "jeanexemple6@outlookcom"
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:jeanexemple6@outlookcom")

# Normalise every hyperlinked cell back onto the shared "Lien hypertexte"
# style (Hyperlinks.Add silently forks a near-duplicate cell style, so pull
# everything back onto the one that already existed in the workbook).
$ws.Range("B2:B7").Style = "Lien hypertexte"

$ws.Range("D11").Select()
